$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: assign directly.
$ws.Range("D2").Value = "42.010.79"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "2.184.59"
$ws.Range("E3").Value = "  -2.81%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -1.80%  "
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("E7").Value = "  -1.95%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -2.95%  "
$ws.Range("E10").Value = "  -5.26%  "
$ws.Range("E11").Value = "  -4.71%  "
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("E13").Value = "  -2.75%  "
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").Value = "2.508.84"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").Value = "2.165.80"
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("E18").Value = "  -6.57%  "
$ws.Range("D19").Value = "41.822.16"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("E21").Value = "  -3.46%  "
$ws.Range("E22").Value = "  -6.58%  "
$ws.Range("E23").Value = "  -9.76%  "
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -6.02%  "
$ws.Range("E28").Value = "  -9.55%  "
$ws.Range("E29").Value = "  -3.36%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("E33").Value = "  +10.56%  "
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("E35").Value = "  -6.20%  "
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("E38").Value = "  -5.72%  "
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("E40").Value = "  -7.54%  "
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("E42").Value = "  -5.66%  "
$ws.Range("E43").Value = "  -8.35%  "
$ws.Range("E44").Value = "  -3.99%  "
$ws.Range("E45").Value = "  -2.94%  "
$ws.Range("E46").Value = "  -3.82%  "
$ws.Range("E47").Value = "  -6.13%  "
$ws.Range("E48").Value = "  -4.81%  "
$ws.Range("E49").Value = "  -4.74%  "
$ws.Range("E50").Value = "  -5.46%  "
$ws.Range("E51").Value = "  +12.95%  "

# Numeric-looking values that must stay stored as text (to match the
# original inline-string cell type rather than becoming a number):
# force text format, assign, then reset style so no stray number format
# style sticks around on the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.586"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0915"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.101"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.787"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000104"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "227.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0781"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.105"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0313"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "59.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.192"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.420"
$ws.Range("D51").Style = "Normal"
